$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume/Change (E) columns with refreshed crypto market data.
# Numeric-looking Price values are protected from Excel auto-number-conversion
# by temporarily forcing a Text format, then the style is restored to Normal so
# the cell keeps its original (unstyled) appearance.

$ws.Range("D2").Value = '25.638.72'
$ws.Range("E2").Value = '  +1.97%  '

$ws.Range("D3").Value = '1.670.38'
$ws.Range("E3").Value = '  +1.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9979'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9989'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4817'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.57%  '

$ws.Range("E8").Value = '  +1.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06168'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07114'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.43%  '

$ws.Range("D11").Value = '1.663.87'
$ws.Range("E11").Value = '  +0.75%  '

$ws.Range("E12").Value = '  +4.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6012'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.423'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.21%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '74.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9989'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9984'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.39%  '

$ws.Range("D18").Value = '25.613.13'
$ws.Range("E18").Value = '  +1.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006811'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.491'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.95%  '

$ws.Range("D22").Value = '1.878.18'
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.731'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.381'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.97%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.405'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '105.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.707'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.989'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.686'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07691'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04376'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9980'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.613'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6205'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9550'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.80%  '

$ws.Range("E38").Value = '  +0.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8680'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9990'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01517'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.876'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.10'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3795'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.95%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.671'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1125'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.259'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05260'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.64%  '

# Rows 50 and 51 swap order: EnergySwap now ranks above Decentraland in the
# table, and both rows receive refreshed Price/Volume figures.
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.416'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.94%  '

$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3363'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.87%  '
